$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns keep their text formatting (avoid Excel auto-numeric conversion)
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '64.866.41'
$ws.Range("E2").Value = '  -0.34%  '
$ws.Range("D3").Value = '3.553.52'
$ws.Range("E3").Value = '  +2.14%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '597.67'
$ws.Range("E5").Value = '  +1.62%  '
$ws.Range("D6").Value = '135.10'
$ws.Range("E6").Value = '  -1.46%  '
$ws.Range("D7").Value = '3.550.37'
$ws.Range("E7").Value = '  +2.07%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("E9").Value = '  +0.50%  '
$ws.Range("E10").Value = '  +0.04%  '
$ws.Range("D11").Value = '6.94'
$ws.Range("E11").Value = '  -2.68%  '
$ws.Range("E12").Value = '  -0.14%  '
$ws.Range("D13").Value = '4.157.53'
$ws.Range("E13").Value = '  +2.14%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '3.565.59'
$ws.Range("E14").Value = '  +2.67%  '
$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D15").Value = '0.0000181'
$ws.Range("E15").Value = '  +0.23%  '
$ws.Range("D16").Value = '26.85'
$ws.Range("E16").Value = '  +1.03%  '
$ws.Range("E17").Value = '  +0.35%  '
$ws.Range("D18").Value = '64.993.75'
$ws.Range("E18").Value = '  -0.09%  '
$ws.Range("D19").Value = '9.93'
$ws.Range("E19").Value = '  +2.12%  '
$ws.Range("D20").Value = '14.31'
$ws.Range("E20").Value = '  +2.72%  '
$ws.Range("D21").Value = '5.79'
$ws.Range("E21").Value = '  +0.50%  '
$ws.Range("D22").Value = '388.03'
$ws.Range("E22").Value = '  -0.18%  '
$ws.Range("E23").Value = '  +3.43%  '
$ws.Range("D24").Value = '3.698.23'
$ws.Range("E24").Value = '  +2.20%  '
$ws.Range("D25").Value = '73.65'
$ws.Range("E25").Value = '  +1.59%  '
$ws.Range("E26").Value = '  +0.08%  '
$ws.Range("E27").Value = '  +3.14%  '
$ws.Range("D28").Value = '7.68'
$ws.Range("E28").Value = '  +3.79%  '
$ws.Range("E29").Value = '  +0.12%  '
$ws.Range("E30").Value = '  +2.76%  '
$ws.Range("E31").Value = '  +2.97%  '
$ws.Range("E32").Value = '  +24.28%  '
$ws.Range("D33").Value = '3.551.33'
$ws.Range("E33").Value = '  +1.61%  '
$ws.Range("E34").Value = '  +0.01%  '
$ws.Range("D35").Value = '23.91'
$ws.Range("E35").Value = '  +3.78%  '
$ws.Range("E36").Value = '  +0.02%  '
$ws.Range("E37").Value = '  -0.95%  '
$ws.Range("D38").Value = '6.88'
$ws.Range("E38").Value = '  +0.63%  '
$ws.Range("D39").Value = '1.54'
$ws.Range("E39").Value = '  +4.36%  '
$ws.Range("D40").Value = '4.95'
$ws.Range("E40").Value = '  +4.81%  '
$ws.Range("D41").Value = '0.0802'
$ws.Range("E41").Value = '  +3.22%  '
$ws.Range("B42").Value = 'EnergySwap'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D42").Value = '26.91'
$ws.Range("E42").Value = '  +8.03%  '
$ws.Range("B43").Value = 'Mantle'
$ws.Range("C43").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D43").Value = '0.822'
$ws.Range("E43").Value = '  +1.42%  '
$ws.Range("D44").Value = '42.67'
$ws.Range("E44").Value = '  +0.05%  '
$ws.Range("E45").Value = '  +0.06%  '
$ws.Range("D46").Value = '4.43'
$ws.Range("E46").Value = '  +1.97%  '
$ws.Range("D47").Value = '1.20'
$ws.Range("E47").Value = '  +3.08%  '
$ws.Range("E48").Value = '  +0.97%  '
$ws.Range("D49").Value = '2.469.72'
$ws.Range("E49").Value = '  +11.47%  '
$ws.Range("D50").Value = '6.88'
$ws.Range("E50").Value = '  +2.77%  '
$ws.Range("D51").Value = '0.864'
$ws.Range("E51").Value = '  +8.05%  '
